$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 1) used "<name>_old" / "<name>_new" suffixes to label
# columns coming from the two compared AHB format versions. Rename them to
# use the concrete format-version identifiers instead ("_FV2210" / "_FV2304").
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the data range into a proper Excel Table ("ListObject") with the
# renamed headers as its column names, plus an autofilter.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U58"))
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (pane split below row 1).
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
